# This workbook holds a per-team state-transition probability matrix
# ("Hampton_A"). Each data row (2-19) represents a "from" state, and the
# columns (B-S) hold the fraction of observed transitions ("games") that
# moved to each "to" state; every row sums to 1.
#
# The commit "added more games, sped up simulate game logic, and drafted
# optimization logic" re-ran the simulation with additional game data,
# which shifted the underlying transition counts and therefore the
# probabilities recorded in this sheet. This script simply writes the
# refreshed probabilities into the matching cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (from state Af0)
$ws.Range("B2").Value = 0.2162162162162162
$ws.Range("C2").Value = 0.528957528957529
$ws.Range("P2").Value = 0.1621621621621622
$ws.Range("S2").Value = 0.09266409266409266

# Row 3 (from state Af1)
$ws.Range("B3").Value = 0.007246376811594203
$ws.Range("C3").Value = 0.007246376811594203
$ws.Range("J3").Value = 0.02173913043478261
$ws.Range("P3").Value = 0.7318840579710145
$ws.Range("S3").Value = 0.2318840579710145

# Row 4 (from state Af2)
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("O4").Value = 0.02631578947368421
$ws.Range("P4").Value = 0.5789473684210527
$ws.Range("S4").Value = 0.3421052631578947

# Row 5 (from state Af3)
$ws.Range("S5").Value = 1

# Row 6 (from state Ai0)
$ws.Range("B6").Value = 0.03763440860215054
$ws.Range("F6").Value = 0.02688172043010753
$ws.Range("J6").Value = 0.3387096774193548
$ws.Range("O6").Value = 0.005376344086021506
$ws.Range("Q6").Value = 0.1827956989247312
$ws.Range("R6").Value = 0.07526881720430108

# Row 7 (from state Ai1)
$ws.Range("B7").Value = 0.0847457627118644
$ws.Range("D7").Value = 0.02259887005649718
$ws.Range("F7").Value = 0.06779661016949153
$ws.Range("J7").Value = 0.1129943502824859
$ws.Range("O7").Value = 0.03389830508474576
$ws.Range("Q7").Value = 0.1412429378531073
$ws.Range("R7").Value = 0.07909604519774012
$ws.Range("S7").Value = 0.4576271186440678

# Row 8 (from state Ai2)
$ws.Range("B8").Value = 0.08658008658008658
$ws.Range("D8").Value = 0.008658008658008658
$ws.Range("F8").Value = 0.05627705627705628
$ws.Range("J8").Value = 0.1385281385281385
$ws.Range("O8").Value = 0.01731601731601732
$ws.Range("Q8").Value = 0.1601731601731602
$ws.Range("R8").Value = 0.1038961038961039
$ws.Range("S8").Value = 0.4285714285714285

# Row 9 (from state Ai3)
$ws.Range("B9").Value = 0.1099476439790576
$ws.Range("D9").Value = 0.01570680628272251
$ws.Range("E9").Value = 0.005235602094240838
$ws.Range("F9").Value = 0.06282722513089005
$ws.Range("J9").Value = 0.1047120418848168
$ws.Range("O9").Value = 0.05235602094240838
$ws.Range("Q9").Value = 0.1465968586387434
$ws.Range("R9").Value = 0.09947643979057591
$ws.Range("S9").Value = 0.4031413612565445

# Row 10 (from state Ar0)
$ws.Range("B10").Value = 0.1040763226366002
$ws.Range("D10").Value = 0.02515177797051171
$ws.Range("F10").Value = 0.05724197745013009
$ws.Range("J10").Value = 0.1196877710320902
$ws.Range("O10").Value = 0.02428447528187338
$ws.Range("Q10").Value = 0.2272333044232437
$ws.Range("R10").Value = 0.07198612315698179
$ws.Range("S10").Value = 0.370338248048569

# Row 11 (from state Bf0)
$ws.Range("G11").Value = 0.164
$ws.Range("J11").Value = 0.07199999999999999
$ws.Range("K11").Value = 0.184
$ws.Range("L11").Value = 0.5679999999999999
$ws.Range("S11").Value = 0.012

# Row 12 (from state Bf1)
$ws.Range("G12").Value = 0.7647058823529411
$ws.Range("J12").Value = 0.1176470588235294
$ws.Range("K12").Value = 0.0196078431372549
$ws.Range("L12").Value = 0.06535947712418301
$ws.Range("S12").Value = 0.03267973856209151

# Row 13 (from state Bf2)
$ws.Range("G13").Value = 0.78125
$ws.Range("J13").Value = 0.21875

# Row 15 (from state Bi0)
$ws.Range("F15").Value = 0.02030456852791878
$ws.Range("H15").Value = 0.1472081218274112
$ws.Range("I15").Value = 0.07106598984771574
$ws.Range("J15").Value = 0.3350253807106599
$ws.Range("K15").Value = 0.05076142131979695
$ws.Range("M15").Value = 0.01522842639593909
$ws.Range("O15").Value = 0.04060913705583756
$ws.Range("S15").Value = 0.3197969543147208

# Row 16 (from state Bi1)
$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.09316770186335403
$ws.Range("J16").Value = 0.453416149068323
$ws.Range("K16").Value = 0.09316770186335403
$ws.Range("M16").Value = 0.02484472049689441
$ws.Range("O16").Value = 0.06832298136645963
$ws.Range("S16").Value = 0.124223602484472

# Row 17 (from state Bi2)
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.1761904761904762
$ws.Range("I17").Value = 0.08333333333333333
$ws.Range("J17").Value = 0.4214285714285714
$ws.Range("K17").Value = 0.09047619047619047
$ws.Range("M17").Value = 0.01904761904761905
$ws.Range("N17").Value = 0.002380952380952381
$ws.Range("O17").Value = 0.05476190476190476
$ws.Range("S17").Value = 0.1238095238095238

# Row 18 (from state Bi3)
$ws.Range("F18").Value = 0.01694915254237288
$ws.Range("H18").Value = 0.1977401129943503
$ws.Range("I18").Value = 0.1016949152542373
$ws.Range("J18").Value = 0.384180790960452
$ws.Range("K18").Value = 0.0847457627118644
$ws.Range("M18").Value = 0.01129943502824859
$ws.Range("N18").Value = 0.005649717514124294
$ws.Range("O18").Value = 0.06779661016949153
$ws.Range("S18").Value = 0.1299435028248588

# Row 19 (from state Br0)
$ws.Range("F19").Value = 0.01885245901639344
$ws.Range("H19").Value = 0.2508196721311475
$ws.Range("I19").Value = 0.08852459016393442
$ws.Range("J19").Value = 0.3524590163934426
$ws.Range("K19").Value = 0.09754098360655737
$ws.Range("M19").Value = 0.0139344262295082
$ws.Range("N19").Value = 0.001639344262295082
$ws.Range("O19").Value = 0.05409836065573771
$ws.Range("S19").Value = 0.1221311475409836
